# Revert "Powerpoint writer: consolidate text run nodes."
#
# Split the trailing-space-bearing runs "Header ", "with ", "Syntax ",
# "Two " and "column " back into separate "word" + "space" runs, matching
# the pre-consolidation run layout. The run properties (<a:rPr/>) on all
# of these runs are empty, so re-setting a sub-range's .Text to its own
# value is enough to make the engine split the run without introducing
# any formatting differences.

$p = $ppt.ActivePresentation

# Slide 1 title: "Header with inline code"
#   "Header " + "with " + "inline code"
#-> "Header" + " " + "with" + " " + "inline code"
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Characters(1, 6).Text = "Header"
$tr1.Characters(7, 1).Text = " "
$tr1.Characters(8, 4).Text = "with"
$tr1.Characters(12, 1).Text = " "

# Slide 2 title: "Syntax highlighting"
#   "Syntax " + "highlighting"
#-> "Syntax" + " " + "highlighting"
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$tr2.Characters(1, 6).Text = "Syntax"
$tr2.Characters(7, 1).Text = " "

# Slide 3 title: "Two column slide"
#   "Two " + "column " + "slide"
#-> "Two" + " " + "column" + " " + "slide"
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Characters(1, 3).Text = "Two"
$tr3.Characters(4, 1).Text = " "
$tr3.Characters(5, 6).Text = "column"
$tr3.Characters(11, 1).Text = " "
